# Auto-generated: update computed market-price / profit columns (H-N)
# on the per-job "Profits" worksheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 245.39131
$ws.Range("I28").Value = 215
$ws.Range("J28").Value = 448
$ws.Range("K28").Value = 215
$ws.Range("L28").Value = 448
$ws.Range("M28").Value = 270
$ws.Range("N28").Value = -1418
$ws.Range("H40").Value = 1418.92
$ws.Range("I40").Value = 1322.75
$ws.Range("K40").Value = 1322.75
$ws.Range("M40").Value = -1147.75
$ws.Range("H86").Value = 2108.3845
$ws.Range("I86").Value = 2140.5
$ws.Range("J86").Value = 2001.3334
$ws.Range("K86").Value = 2140.5
$ws.Range("L86").Value = 2001.3334
$ws.Range("M86").Value = -1017.5
$ws.Range("N86").Value = -4247.3334
$ws.Range("H89").Value = 2108.3845
$ws.Range("I89").Value = 2140.5
$ws.Range("J89").Value = 2001.3334
$ws.Range("K89").Value = 10702.5
$ws.Range("L89").Value = 10006.667
$ws.Range("M89").Value = -5086.5
$ws.Range("N89").Value = -21238.667
$ws.Range("H137").Value = 1497.75
$ws.Range("I137").Value = 1428.5927
$ws.Range("J137").Value = 1607.5883
$ws.Range("K137").Value = 4285.7781
$ws.Range("L137").Value = 4822.7649
$ws.Range("M137").Value = -1735.7781
$ws.Range("N137").Value = -9922.7649
$ws.Range("H138").Value = 2325.66
$ws.Range("I138").Value = 1262.1666
$ws.Range("J138").Value = 3920.9
$ws.Range("K138").Value = 3786.4998
$ws.Range("L138").Value = 11762.7
$ws.Range("M138").Value = 1353.5002
$ws.Range("N138").Value = -22042.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 26288.75
$ws.Range("J24").Value = 26288.75
$ws.Range("L24").Value = 26288.75
$ws.Range("N24").Value = -27036.75
$ws.Range("H32").Value = 974.6799999999999
$ws.Range("I32").Value = 887.875
$ws.Range("J32").Value = 1611.25
$ws.Range("K32").Value = 887.875
$ws.Range("L32").Value = 1611.25
$ws.Range("M32").Value = -600.875
$ws.Range("N32").Value = -2185.25
$ws.Range("H63").Value = 4090.9
$ws.Range("I63").Value = 4444.143
$ws.Range("K63").Value = 4444.143
$ws.Range("M63").Value = -3758.143
$ws.Range("H66").Value = 4090.9
$ws.Range("I66").Value = 4444.143
$ws.Range("K66").Value = 22220.715
$ws.Range("M66").Value = -18788.715
$ws.Range("H74").Value = 998.1786
$ws.Range("I74").Value = 941.96
$ws.Range("K74").Value = 941.96
$ws.Range("M74").Value = -67.96000000000004
$ws.Range("H77").Value = 998.1786
$ws.Range("I77").Value = 941.96
$ws.Range("K77").Value = 4709.8
$ws.Range("M77").Value = -341.8000000000002
$ws.Range("H100").Value = 26288.75
$ws.Range("J100").Value = 26288.75
$ws.Range("L100").Value = 26288.75
$ws.Range("N100").Value = -28452.75
$ws.Range("H122").Value = 1528.6
$ws.Range("I122").Value = 1380.1666
$ws.Range("J122").Value = 1751.25
$ws.Range("K122").Value = 4140.4998
$ws.Range("L122").Value = 5253.75
$ws.Range("M122").Value = -1690.4998
$ws.Range("N122").Value = -10153.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 135173.8
$ws.Range("I86").Value = 2263
$ws.Range("J86").Value = 500678.5
$ws.Range("K86").Value = 2263
$ws.Range("L86").Value = 500678.5
$ws.Range("M86").Value = -1140
$ws.Range("N86").Value = -502924.5
$ws.Range("H89").Value = 135173.8
$ws.Range("I89").Value = 2263
$ws.Range("J89").Value = 500678.5
$ws.Range("K89").Value = 11315
$ws.Range("L89").Value = 2503392.5
$ws.Range("M89").Value = -5699
$ws.Range("N89").Value = -2514624.5
$ws.Range("H94").Value = 73556.86
$ws.Range("I94").Value = 1779.6
$ws.Range("J94").Value = 253000
$ws.Range("K94").Value = 1779.6
$ws.Range("L94").Value = 253000
$ws.Range("M94").Value = -1328.6
$ws.Range("N94").Value = -253902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1687.8871
$ws.Range("I31").Value = 1220.4634
$ws.Range("K31").Value = 1220.4634
$ws.Range("M31").Value = -925.4634000000001
$ws.Range("H34").Value = 1687.8871
$ws.Range("I34").Value = 1220.4634
$ws.Range("K34").Value = 1220.4634
$ws.Range("M34").Value = -1018.4634
$ws.Range("H58").Value = 1483432.6
$ws.Range("I58").Value = 3368721.8
$ws.Range("K58").Value = 3368721.8
$ws.Range("M58").Value = -3368518.8
$ws.Range("H60").Value = 31720
$ws.Range("I60").Value = 3000
$ws.Range("K60").Value = 3000
$ws.Range("M60").Value = -2489
$ws.Range("H62").Value = 57906.11
$ws.Range("I62").Value = 73729.28999999999
$ws.Range("J62").Value = 2525
$ws.Range("K62").Value = 73729.28999999999
$ws.Range("L62").Value = 2525
$ws.Range("M62").Value = -73105.28999999999
$ws.Range("N62").Value = -3773
$ws.Range("H65").Value = 57906.11
$ws.Range("I65").Value = 73729.28999999999
$ws.Range("J65").Value = 2525
$ws.Range("K65").Value = 368646.45
$ws.Range("L65").Value = 12625
$ws.Range("M65").Value = -365526.45
$ws.Range("N65").Value = -18865
$ws.Range("H132").Value = 323022.5
$ws.Range("I132").Value = 398444.97
$ws.Range("K132").Value = 1195334.91
$ws.Range("M132").Value = -1192804.91
$ws.Range("H136").Value = 1483432.6
$ws.Range("I136").Value = 3368721.8
$ws.Range("K136").Value = 10106165.4
$ws.Range("M136").Value = -10103615.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 46347.91
$ws.Range("I117").Value = 643.6
$ws.Range("J117").Value = 84434.836
$ws.Range("K117").Value = 1930.8
$ws.Range("L117").Value = 253304.508
$ws.Range("M117").Value = 1511.2
$ws.Range("N117").Value = -260188.508

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2521.077
$ws.Range("I16").Value = 2570.182
$ws.Range("J16").Value = 2251
$ws.Range("K16").Value = 2570.182
$ws.Range("L16").Value = 2251
$ws.Range("M16").Value = -2400.182
$ws.Range("N16").Value = -2591
$ws.Range("H46").Value = 1641.5
$ws.Range("I46").Value = 1414
$ws.Range("J46").Value = 1960
$ws.Range("K46").Value = 1414
$ws.Range("L46").Value = 1960
$ws.Range("M46").Value = -1226
$ws.Range("N46").Value = -2336
$ws.Range("H68").Value = 2608.3845
$ws.Range("I68").Value = 1625.75
$ws.Range("J68").Value = 4180.6
$ws.Range("K68").Value = 1625.75
$ws.Range("L68").Value = 4180.6
$ws.Range("M68").Value = -876.75
$ws.Range("N68").Value = -5678.6
$ws.Range("H71").Value = 2608.3845
$ws.Range("I71").Value = 1625.75
$ws.Range("J71").Value = 4180.6
$ws.Range("K71").Value = 8128.75
$ws.Range("L71").Value = 20903
$ws.Range("M71").Value = -4384.75
$ws.Range("N71").Value = -28391
$ws.Range("H132").Value = 4543.115
$ws.Range("I132").Value = 4407.2
$ws.Range("K132").Value = 13221.6
$ws.Range("M132").Value = -10691.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 32505.25
$ws.Range("J15").Value = 32505.25
$ws.Range("L15").Value = 32505.25
$ws.Range("N15").Value = -33081.25
$ws.Range("H107").Value = 549.9583
$ws.Range("I107").Value = 505.05264
$ws.Range("J107").Value = 720.6
$ws.Range("K107").Value = 1515.15792
$ws.Range("L107").Value = 2161.8
$ws.Range("M107").Value = 404.8420799999999
$ws.Range("N107").Value = -6001.8
$ws.Range("H125").Value = 60715
$ws.Range("J125").Value = 60715
$ws.Range("L125").Value = 60715
$ws.Range("N125").Value = -70555
